$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-based values for rows 2-7 (Sost -> Lrp6 pairs)
$ws.Range("G2").Value = 0.1212753333333333
$ws.Range("H2").Value = 0.363826
$ws.Range("M2").Value = 12.80871533333333
$ws.Range("N2").Value = 38.426146
$ws.Range("O2").Value = 0.1716721220213608
$ws.Range("P2").Value = 0.1716721220213608
$ws.Range("Q2").Value = 1.553381221621778
$ws.Range("R2").Value = 13.980430994596
$ws.Range("S2").Value = 0.1716721220213608
$ws.Range("T2").Value = 0.1716721220213608

$ws.Range("G3").Value = 0.1212753333333333
$ws.Range("H3").Value = 0.363826
$ws.Range("M3").Value = 21.85073566666667
$ws.Range("N3").Value = 65.552207
$ws.Range("O3").Value = 0.2928601395225403
$ws.Range("P3").Value = 0.2928601395225403
$ws.Range("Q3").Value = 2.649955251553555
$ws.Range("R3").Value = 23.849597263982
$ws.Range("S3").Value = 0.2928601395225403
$ws.Range("T3").Value = 0.2928601395225403

$ws.Range("G4").Value = 0.1212753333333333
$ws.Range("H4").Value = 0.363826
$ws.Range("M4").Value = 10.846871
$ws.Range("N4").Value = 32.540613
$ws.Range("O4").Value = 0.1453779956383313
$ws.Range("P4").Value = 0.1453779956383313
$ws.Range("Q4").Value = 1.315457896148667
$ws.Range("R4").Value = 11.839121065338
$ws.Range("S4").Value = 0.1453779956383313
$ws.Range("T4").Value = 0.1453779956383313

$ws.Range("G5").Value = 0.1212753333333333
$ws.Range("H5").Value = 0.363826
$ws.Range("M5").Value = 7.389532
$ws.Range("N5").Value = 22.168596
$ws.Range("O5").Value = 0.09904011496636306
$ws.Range("P5").Value = 0.09904011496636304
$ws.Range("Q5").Value = 0.8961679564773333
$ws.Range("R5").Value = 8.065511608295999
$ws.Range("S5").Value = 0.09904011496636306
$ws.Range("T5").Value = 0.09904011496636304

$ws.Range("G6").Value = 0.1212753333333333
$ws.Range("H6").Value = 0.363826
$ws.Range("M6").Value = 7.190038666666666
$ws.Range("N6").Value = 21.570116
$ws.Range("O6").Value = 0.09636635393950015
$ws.Range("P6").Value = 0.09636635393950013
$ws.Range("Q6").Value = 0.8719743359795555
$ws.Range("R6").Value = 7.847769023815999
$ws.Range("S6").Value = 0.09636635393950015
$ws.Range("T6").Value = 0.09636635393950013

$ws.Range("G7").Value = 0.1212753333333333
$ws.Range("H7").Value = 0.363826
$ws.Range("M7").Value = 14.525612
$ws.Range("N7").Value = 43.576836
$ws.Range("O7").Value = 0.1946832739119044
$ws.Range("P7").Value = 0.1946832739119044
$ws.Range("Q7").Value = 1.761598437170667
$ws.Range("R7").Value = 15.854385934536
$ws.Range("S7").Value = 0.1946832739119044
$ws.Range("T7").Value = 0.1946832739119044
